$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 6135
$ws.Cells.Item(3, 12).Value = 6667
$ws.Cells.Item(4, 12).Value = 1641
$ws.Cells.Item(6, 12).Value = 5472
$ws.Cells.Item(7, 12).Value = 20312

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(5, 12).Value = 72
$ws.Cells.Item(8, 12).Value = 1345
$ws.Cells.Item(9, 12).Value = 118
$ws.Cells.Item(16, 12).Value = 46
$ws.Cells.Item(18, 12).Value = 138
$ws.Cells.Item(29, 12).Value = 1138
$ws.Cells.Item(33, 12).Value = 917
$ws.Cells.Item(36, 12).Value = 260
$ws.Cells.Item(37, 12).Value = 779
$ws.Cells.Item(42, 12).Value = 647
$ws.Cells.Item(43, 12).Value = 152
$ws.Cells.Item(44, 12).Value = 138
$ws.Cells.Item(52, 12).Value = 433
$ws.Cells.Item(53, 12).Value = 225
$ws.Cells.Item(63, 12).Value = 61
$ws.Cells.Item(67, 12).Value = 705
$ws.Cells.Item(68, 12).Value = 65
$ws.Cells.Item(76, 12).Value = 308
$ws.Cells.Item(78, 12).Value = 267
$ws.Cells.Item(79, 12).Value = 561
$ws.Cells.Item(83, 12).Value = 445
$ws.Cells.Item(85, 12).Value = 1009
$ws.Cells.Item(86, 12).Value = 131
$ws.Cells.Item(91, 12).Value = 275
$ws.Cells.Item(95, 12).Value = 288
$ws.Cells.Item(96, 12).Value = 226
$ws.Cells.Item(98, 12).Value = 109
$ws.Cells.Item(99, 12).Value = 354
$ws.Cells.Item(100, 12).Value = 39
$ws.Cells.Item(101, 12).Value = 20312

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 12).Value = 72
$ws.Cells.Item(7, 12).Value = 226

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 12).Value = 306
$ws.Cells.Item(3, 12).Value = 419
$ws.Cells.Item(7, 12).Value = 1009

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(6, 12).Value = 123
$ws.Cells.Item(7, 12).Value = 433

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 12).Value = 76
$ws.Cells.Item(7, 12).Value = 225

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 12).Value = 404
$ws.Cells.Item(3, 12).Value = 477
$ws.Cells.Item(7, 12).Value = 1345

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(6, 12).Value = 99
$ws.Cells.Item(7, 12).Value = 445

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 12).Value = 248
$ws.Cells.Item(7, 12).Value = 917

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(3, 12).Value = 93
$ws.Cells.Item(7, 12).Value = 288

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 12).Value = 235
$ws.Cells.Item(7, 12).Value = 779

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 12).Value = 144
$ws.Cells.Item(7, 12).Value = 354

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 12).Value = 200
$ws.Cells.Item(6, 12).Value = 163
$ws.Cells.Item(7, 12).Value = 705

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 12).Value = 440
$ws.Cells.Item(6, 12).Value = 279
$ws.Cells.Item(7, 12).Value = 1138

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 12).Value = 34
$ws.Cells.Item(7, 12).Value = 138

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 12).Value = 137
$ws.Cells.Item(7, 12).Value = 308

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 12).Value = 177
$ws.Cells.Item(6, 12).Value = 180
$ws.Cells.Item(7, 12).Value = 647

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(6, 12).Value = 76
$ws.Cells.Item(7, 12).Value = 267

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(2, 12).Value = 94
$ws.Cells.Item(3, 12).Value = 125
$ws.Cells.Item(4, 12).Value = 15
$ws.Cells.Item(7, 12).Value = 275

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 12).Value = 180
$ws.Cells.Item(7, 12).Value = 561

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(2, 12).Value = 49
$ws.Cells.Item(7, 12).Value = 138

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 12).Value = 84
$ws.Cells.Item(7, 12).Value = 260

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Cells.Item(6, 12).Value = 24
$ws.Cells.Item(7, 12).Value = 39

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(6, 12).Value = 49
$ws.Cells.Item(7, 12).Value = 109

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(6, 12).Value = 28
$ws.Cells.Item(7, 12).Value = 118

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(6, 12).Value = 32
$ws.Cells.Item(7, 12).Value = 72

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(4, 12).Value = 70
$ws.Cells.Item(7, 12).Value = 131

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 12).Value = 70
$ws.Cells.Item(6, 12).Value = 61

$ws = $wb.Worksheets.Item('North Park')
$ws.Cells.Item(2, 12).Value = 22
$ws.Cells.Item(7, 12).Value = 65

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 12).Value = 47
$ws.Cells.Item(7, 12).Value = 152

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(6, 12).Value = 30
$ws.Cells.Item(7, 12).Value = 46
